# Add two new localization entries ("tag_conductive" / "tag_non_conductive")
# right after the existing "tag_sink" row, pushing the material_* rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 12 (material_clay),
# shifting all rows from 12 downward by two positions.
$ws.Rows("12:13").Insert()

# Fill in the two new key/value pairs in the freshly inserted rows.
$ws.Range("A12").Value = "tag_conductive"
$ws.Range("B12").Value = "Conductive"

$ws.Range("A13").Value = "tag_non_conductive"
$ws.Range("B13").Value = "Non-Conductive"

# Leave the selection on the last edited cell, like the author did.
$ws.Range("B13").Select()
